$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: wrap "Clone Phishing " in a bookmark (_Hlk184124852), splitting
# the existing run into two runs around the bookmark.
# ---------------------------------------------------------------------------
$clonePara = $d.Paragraphs(6)
$cloneRange = $clonePara.Range
$cloneRange.Find.Execute("Clone Phishing ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_Hlk184124852", $cloneRange)

# ---------------------------------------------------------------------------
# Change 2: after the "Simulation:" paragraph, add a blank bold paragraph, a
# bold "TO DO:" paragraph, and six new bulleted (numId 4) "to do" items.
# ---------------------------------------------------------------------------
$simulationPara = $d.Paragraphs(15)

# A list paragraph that already uses the numId=4 bulleted list, used below as
# a formatting template so the new bullet items pick up the exact same
# pPr/numPr/rPr pattern as the existing list.
$listTemplate = $d.Paragraphs(14)

# Blank bold paragraph right after "Simulation:"
$simulationPara.Range.InsertParagraphAfter()

# "TO DO:" bold paragraph right after the blank paragraph
$blankPara = $d.Paragraphs(16)
$blankPara.Range.InsertParagraphAfter()
$todoPara = $d.Paragraphs(17)
$todoPara.Range.Text = "TO DO:"

$items = @(
    "Complete welcome instructions",
    "Show example of phishing email",
    "Add 3 real phishing emails",
    "Add choices at the bottom to show what kind of phishing email it is.",
    "Add scoring system",
    "Add social engineering email"
)

$afterIndex = 17
foreach ($item in $items) {
    $listTemplate.Range.Copy()
    $nextPara = $d.Paragraphs($afterIndex + 1)
    $insertPoint = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
    $insertPoint.Paste()

    $newPara = $d.Paragraphs($afterIndex + 1)
    $textRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $textRange.Text = $item

    $afterIndex = $afterIndex + 1
}
